$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.28"
$ws.Range("E2").Value = "'-1.18%"
$ws.Range("D3").Value = "'35.86"
$ws.Range("E3").Value = "'-5.00%"
$ws.Range("D4").Value = "'5.114"
$ws.Range("E4").Value = "'-0.52%"
$ws.Range("D5").Value = "'0.07709"
$ws.Range("E5").Value = "'-2.46%"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").Value = "'8.288"
$ws.Range("E6").Value = "'0.10%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.854"
$ws.Range("E7").Value = "'-2.89%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.942"
$ws.Range("E8").Value = "'-5.01%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9205"
$ws.Range("E9").Value = "'-0.43%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1115"
$ws.Range("E10").Value = "'-7.55%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1855"
$ws.Range("E11").Value = "'-3.97%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08762"
$ws.Range("E12").Value = "'-4.28%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03322"
$ws.Range("E13").Value = "'0.58%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09542"
$ws.Range("E14").Value = "'-0.92%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001388"
$ws.Range("E15").Value = "'0.71%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006121"
$ws.Range("E16").Value = "'4.35%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.365"
$ws.Range("E17").Value = "'-4.42%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.387"
$ws.Range("E18").Value = "'-0.60%"
$ws.Range("D19").Value = "'0.3445"
$ws.Range("E19").Value = "'1.32%"
$ws.Range("D20").Value = "'6.312"
$ws.Range("E20").Value = "'19.69%"
$ws.Range("D21").Value = "'0.1298"
$ws.Range("E21").Value = "'2.00%"
$ws.Range("D22").Value = "'0.2312"
$ws.Range("E22").Value = "'-10.69%"
$ws.Range("D23").Value = "'0.04336"
$ws.Range("E23").Value = "'-0.54%"
$ws.Range("D24").Value = "'0.001204"
$ws.Range("E24").Value = "'-3.48%"
$ws.Range("D25").Value = "'0.004258"
$ws.Range("E25").Value = "'-1.25%"
$ws.Range("D26").Value = "'0.0001333"
$ws.Range("E26").Value = "'9.35%"
$ws.Range("D27").Value = "'0.0002906"
$ws.Range("D39").Value = "'0.02095"
$ws.Range("E39").Value = "'-1.39%"
$ws.Range("D40").Value = "'0.04917"
$ws.Range("D41").Value = "'0.007538"
$ws.Range("E41").Value = "'-1.06%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'-1.04%"
$ws.Range("D43").Value = "'0.008552"
$ws.Range("E43").Value = "'-5.83%"
$ws.Range("D44").Value = "'0.002074"
$ws.Range("E44").Value = "'3.30%"
$ws.Range("D45").Value = "'0.008409"
$ws.Range("E45").Value = "'-2.17%"
$ws.Range("D46").Value = "'0.00006454"
$ws.Range("E46").Value = "'-3.66%"
$ws.Range("E47").Value = "'0.29%"
$ws.Range("D48").Value = "'0.003302"
$ws.Range("E48").Value = "'15.36%"
$ws.Range("D49").Value = "'0.001446"
$ws.Range("E49").Value = "'20.61%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.29%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.29%"
